$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Edit 1: paragraph 4 - "is available " -> "is availab" + _GoBack + "le "
#   (the _GoBack bookmark moves from the end of the paragraph to here)
# ------------------------------------------------------------------
$p4 = $d.Paragraphs.Item(4)
$p4text = $p4.Range.Text
$p4start = $p4.Range.Start
$availIdx = $p4text.IndexOf("is available")
$availStart = $p4start + $availIdx

# split "is available " into "is availab" | "le " without altering the text,
# using a no-op formatting toggle so the run boundary is recomputed cleanly
$splitPoint = $availStart + 10
$rSplit = $d.Range($availStart, $splitPoint)
$rSplit.Font.Bold = 1
$rSplit.Font.Bold = 0

# move the hidden _GoBack bookmark to the new split point (this both creates
# it here and removes it from wherever it previously was)
$rPoint = $d.Range($splitPoint, $splitPoint)
$d.Bookmarks.Add("_GoBack", $rPoint)

# ------------------------------------------------------------------
# Edit 2: paragraph 34 - " services" -> " Services" (capitalise the S),
#   split into " S" | "ervices"
# ------------------------------------------------------------------
$p34 = $d.Paragraphs.Item(34)
$p34start = $p34.Range.Start
$p34text = $p34.Range.Text
$svcIdx = $p34text.IndexOf(" services")
$svcStart = $p34start + $svcIdx

# fix the capitalisation first (this is a real content edit and will merge
# the paragraph's same-formatted runs together as a side effect)
$sStart = $svcStart + 1
$sEnd = $sStart + 1
$rChar = $d.Range($sStart, $sEnd)
$rChar.Text = "S"

# now restore/produce the desired run boundaries with no-op formatting
# toggles (these do not touch the text, so they don't cause further merges)
$p34text2 = $p34.Range.Text
$apEnd = $p34start + ($p34text2.IndexOf("Access P") + 8)
$ointEnd = $apEnd + 4
$sEnd2 = $ointEnd + 2

$rB1 = $d.Range($p34start, $apEnd)
$rB1.Font.Bold = 1
$rB1.Font.Bold = 0

$rB2 = $d.Range($apEnd, $ointEnd)
$rB2.Font.Bold = 1
$rB2.Font.Bold = 0

$rB3 = $d.Range($ointEnd, $sEnd2)
$rB3.Font.Bold = 1
$rB3.Font.Bold = 0

# ------------------------------------------------------------------
# Edit 3: paragraph 35 - "Subject Meta-data Publishing" ->
#   "Service Metadata Publishing", split into "Service Meta" | "data Publishing"
# ------------------------------------------------------------------
$p35 = $d.Paragraphs.Item(35)
$p35start = $p35.Range.Start
$p35text = $p35.Range.Text
$subjIdx = $p35text.IndexOf("Subject Meta-data Publishing")
$subjStart = $p35start + $subjIdx
$subjEnd = $subjStart + 29
$rSubj = $d.Range($subjStart, $subjEnd)
$rSubj.Text = "Service Metadata Publishing"

# restore the " - " | "Service Meta" | "data Publishing" run split
$p35text2 = $p35.Range.Text
$dashEnd = $p35start + ($p35text2.IndexOf(" - ") + 3)
$metaEnd = $p35start + ($p35text2.IndexOf("Service Meta") + 12)

$rC1 = $d.Range($p35start, $dashEnd)
$rC1.Font.Bold = 1
$rC1.Font.Bold = 0

$rC2 = $d.Range($dashEnd, $metaEnd)
$rC2.Font.Bold = 1
$rC2.Font.Bold = 0

# ------------------------------------------------------------------
# Edit 4: remove the "DRAFT" watermark shape from the header, leaving an
#   empty paragraph styled "Header"
# ------------------------------------------------------------------
$sec = $d.Sections.Item(1)
$hdr = $sec.Headers.Item(1)
if ($hdr.Shapes.Count -gt 0) {
    for ($i = $hdr.Shapes.Count; $i -ge 1; $i--) {
        $hdr.Shapes.Item($i).Delete()
    }
}
$hp1 = $hdr.Range.Paragraphs.Item(1)
$hp1.Style = "Header"
